# Populate the Config sheet's example row with the sample "weather_data"
# workbook config (sheet name "Sheet1", date column index 0) instead of the
# previous "temperaturi" / 2022 / 184 placeholder values.
#
# Shared-string order matters for a minimal diff against the target file, so
# write column B ("Sheet1") before column A ("weather_data") - that way the
# new unique strings land in the sheet in the same order they appear in the
# target workbook's sharedStrings table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sheet1"
$ws.Range("A2").Value = "weather_data"
$ws.Range("C2").Value = 0

# Move/restore the visible selection to A5, matching the saved view state.
$null = $ws.Range("A5").Select()
